$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Resize existing / new columns (ColumnWidth excludes Excel's fixed padding of 5/6
# character, so subtract it to land exactly on the target "width" stored in the XML).
$pad = 5/6
$ws.Columns.Item(1).ColumnWidth = 30 - $pad
$ws.Columns.Item(2).ColumnWidth = 15 - $pad
$ws.Columns.Item(3).ColumnWidth = 40 - $pad
$ws.Columns.Item(4).ColumnWidth = 25 - $pad
$ws.Columns.Item(5).ColumnWidth = 25 - $pad
$ws.Columns.Item(6).ColumnWidth = 15 - $pad
$ws.Columns.Item(7).ColumnWidth = 15 - $pad
$ws.Columns.Item(8).ColumnWidth = 80 - $pad

# Row 2 - Fernando's Laboratory
$ws.Range("A2").Value = "Fernando's Laboratory"
$ws.Range("B2").Value = "Lab"
$ws.Range("C2").Value = "7 W Lane, Central NY NY 11723"
$ws.Range("D2").Value = "wwww.tcgdex.net"
$ws.Range("E2").Value = "fernhean@hotmail.nw"
$ws.Range("F2").Value = "123-456-7890"
$ws.Range("G2").Value = "312-312-4212"
$ws.Range("H2").Value = "My lab is full of enourmous surprises"

# Row 3 - Tiffy's Bunnies
$ws.Range("A3").Value = "Tiffy's Bunnies"
$ws.Range("B3").Value = "Animal Care"
$ws.Range("C3").Value = "Earth"
$ws.Range("D3").Value = "tiffysbunnies.com"
$ws.Range("E3").Value = "tiffy23@aol.321"
$ws.Range("F3").Value = "653-123-4632"
$ws.Range("G3").Value = "N/A"
$ws.Range("H3").Value = "The most epic bunnies in existance!"
